$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before the current row 124, pushing existing
# rows 124-184 down to 128-188. Excel copies formatting from the row
# above into the newly inserted rows (so column D keeps its date style).
$ws.Range("A124:A127").EntireRow.Insert()

# --- Row 124 ---
$ws.Cells.Item(124, 1).Value = 10
$ws.Cells.Item(124, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(124, 3).Value = "La Araucanía"
$ws.Cells.Item(124, 4).Value = 44518
$ws.Cells.Item(124, 5).Value = 9
$ws.Cells.Item(124, 6).Value = "Fruta"
$ws.Cells.Item(124, 7).Value = 100101
$ws.Cells.Item(124, 8).Value = "Berries"
$ws.Cells.Item(124, 9).Value = 100112025
$ws.Cells.Item(124, 10).Value = "Frutilla"
$ws.Cells.Item(124, 11).Value = "Sin especificar"
$ws.Cells.Item(124, 12).Value = "Primera"
$ws.Cells.Item(124, 13).Value = 5000
$ws.Cells.Item(124, 14).Value = 8000
$ws.Cells.Item(124, 15).Value = 9000
$ws.Cells.Item(124, 16).Value = 8500
$ws.Cells.Item(124, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(124, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(124, 19).Value = 1214
$ws.Cells.Item(124, 20).Value = 7

# --- Row 125 ---
$ws.Cells.Item(125, 1).Value = 10
$ws.Cells.Item(125, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(125, 3).Value = "La Araucanía"
$ws.Cells.Item(125, 4).Value = 44518
$ws.Cells.Item(125, 5).Value = 9
$ws.Cells.Item(125, 6).Value = "Fruta"
$ws.Cells.Item(125, 7).Value = 100101
$ws.Cells.Item(125, 8).Value = "Berries"
$ws.Cells.Item(125, 9).Value = 100112025
$ws.Cells.Item(125, 10).Value = "Frutilla"
$ws.Cells.Item(125, 11).Value = "Sin especificar"
$ws.Cells.Item(125, 12).Value = "Primera"
$ws.Cells.Item(125, 13).Value = 300
$ws.Cells.Item(125, 14).Value = 8000
$ws.Cells.Item(125, 15).Value = 8000
$ws.Cells.Item(125, 16).Value = 8000
$ws.Cells.Item(125, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(125, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(125, 19).Value = 1143
$ws.Cells.Item(125, 20).Value = 7

# --- Row 126 ---
$ws.Cells.Item(126, 1).Value = 10
$ws.Cells.Item(126, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(126, 3).Value = "La Araucanía"
$ws.Cells.Item(126, 4).Value = 44518
$ws.Cells.Item(126, 5).Value = 9
$ws.Cells.Item(126, 6).Value = "Fruta"
$ws.Cells.Item(126, 7).Value = 100101
$ws.Cells.Item(126, 8).Value = "Berries"
$ws.Cells.Item(126, 9).Value = 100112025
$ws.Cells.Item(126, 10).Value = "Frutilla"
$ws.Cells.Item(126, 11).Value = "Sin especificar"
$ws.Cells.Item(126, 12).Value = "Segunda"
$ws.Cells.Item(126, 13).Value = 400
$ws.Cells.Item(126, 14).Value = 7000
$ws.Cells.Item(126, 15).Value = 7500
$ws.Cells.Item(126, 16).Value = 7250
$ws.Cells.Item(126, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(126, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(126, 19).Value = 1036
$ws.Cells.Item(126, 20).Value = 7

# --- Row 127 ---
$ws.Cells.Item(127, 1).Value = 10
$ws.Cells.Item(127, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(127, 3).Value = "La Araucanía"
$ws.Cells.Item(127, 4).Value = 44518
$ws.Cells.Item(127, 5).Value = 9
$ws.Cells.Item(127, 6).Value = "Fruta"
$ws.Cells.Item(127, 7).Value = 100101
$ws.Cells.Item(127, 8).Value = "Berries"
$ws.Cells.Item(127, 9).Value = 100112025
$ws.Cells.Item(127, 10).Value = "Frutilla"
$ws.Cells.Item(127, 11).Value = "Sin especificar"
$ws.Cells.Item(127, 12).Value = "Tercera"
$ws.Cells.Item(127, 13).Value = 100
$ws.Cells.Item(127, 14).Value = 5000
$ws.Cells.Item(127, 15).Value = 5000
$ws.Cells.Item(127, 16).Value = 5000
$ws.Cells.Item(127, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(127, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(127, 19).Value = 714
$ws.Cells.Item(127, 20).Value = 7
